$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New epoch header columns (R1:U1)
$ws.Range("R1").Value = "epoch850"
$ws.Range("S1").Value = "epoch900"
$ws.Range("T1").Value = "epoch950"
$ws.Range("U1").Value = "epoch1000"

# Update existing row 2 values (B2:Q2) to new accuracy value
$newValue = 78.38541641831398
$ws.Range("B2:Q2").Value = $newValue

# New row 2 cells for the added columns (R2:U2)
$ws.Range("R2:U2").Value = $newValue
